$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows above the old header row (old row 3), pushing the
# header + data rows down to rows 7-10.
$ws.Rows("3:6").Insert()

# Fill in the new box-dimension rows that were inserted.
$ws.Range("A3").Value = "Lx="
$ws.Range("B3").Value = 23.41
$ws.Range("A4").Value = "Ly="
$ws.Range("B4").Value = 18.69
$ws.Range("A5").Value = "Lz="
$ws.Range("B5").Value = 14.5

# Remove the now-unused "wj" column (column G) entirely; this shifts the
# "dij" column (old H) left into G for the header and all data rows.
$ws.Columns("G").Delete()

# Update the active selection to match the authored workbook.
[void]$ws.Range("A3:B5").Select()
